# Updated the times of inputs in the cover page
$d = $word.ActiveDocument

# --- 1. Split the title run "ECM1410 Cover page " into three runs, with
#        "page" isolated between a pair of grammar-checker <w:proofErr/>
#        markers (the run split Word leaves behind after re-checking the
#        title on save), without altering the visible text.
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$insertRange = $d.Range($titleRange.Start, $titleRange.End - 1)
$titleXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:jc w:val="center"/></w:pPr><w:r><w:t xml:space="preserve">ECM1410 Cover </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>page</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document>
'@
$insertRange.InsertXML($titleXml)

# --- 2. Fill in the new practical session row (date / time / duration /
#        observer / driver) in the second table.
$t = $d.Tables.Item(2)
$t.Cell(7, 1).Range.Text = "03/3/23"
$t.Cell(7, 2).Range.Text = "13:00"
$t.Cell(7, 3).Range.Text = "1h20"
$t.Cell(7, 4).Range.Text = "Observer"
$t.Cell(7, 5).Range.Text = "Driver"

Write-Output "done"
